# Insert a new price-record row at row 135 on the single data sheet.
# This shifts the existing rows 135:218 down to 136:219 (dimension grows
# from A1:R218 to A1:R219) and fills the newly opened row 135 with the
# new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("135:135").Insert()

$ws.Range("A135").Value2 = 9
$ws.Range("B135").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C135").Value2 = "Metropolitana"
$ws.Range("D135").Value2 = 44719
$ws.Range("E135").Value2 = 13
$ws.Range("F135").Value2 = 100112026
$ws.Range("G135").Value2 = "Haba"
$ws.Range("H135").Value2 = "Sin especificar"
$ws.Range("I135").Value2 = "Primera"
$ws.Range("J135").Value2 = 61
$ws.Range("K135").Value2 = 17000
$ws.Range("L135").Value2 = 18000
$ws.Range("M135").Value2 = 17508
$ws.Range("N135").Value2 = "`$/saco 25 kilos"
$ws.Range("O135").Value2 = "Provincia de Limarí"
$ws.Range("P135").Value2 = 700
$ws.Range("Q135").Value2 = 25
$ws.Range("R135").Value2 = "Hortaliza"
